$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.09132603783370996
$ws.Range("D2").Value = 0.01991732608433239
$ws.Range("E2").Value = 0.1662197965670416
$ws.Range("F2").Value = 0.4407179464367204
$ws.Range("G2").Value = 0.2871559136283182
$ws.Range("H2").Value = 0.4516220209577995
$ws.Range("I2").Value = 0.3980573076856082
$ws.Range("K2").Value = 0.474974610147143
$ws.Range("M2").Value = 0.2413221055721166
$ws.Range("N2").Value = 1.174092548864579
$ws.Range("O2").Value = 1.402229350648625
$ws.Range("B3").Value = 0.08106439087796957
$ws.Range("D3").Value = 0.0175912346774183
$ws.Range("E3").Value = 0.1576725710641469
$ws.Range("F3").Value = 0.4364268309321204
$ws.Range("G3").Value = 0.284227587381892
$ws.Range("H3").Value = 0.4534663577912283
$ws.Range("I3").Value = 0.4036917249837471
$ws.Range("K3").Value = 0.4174331193433716
$ws.Range("M3").Value = 0.2136605678512424
$ws.Range("N3").Value = 1.184789590088933
$ws.Range("O3").Value = 1.399645377150947
$ws.Range("B4").Value = 0.0747639670415623
$ws.Range("D4").Value = 0.01615502861609741
$ws.Range("E4").Value = 0.1525582784427328
$ws.Range("F4").Value = 0.4340943183082615
$ws.Range("G4").Value = 0.2826612807532101
$ws.Range("H4").Value = 0.454815571434807
$ws.Range("I4").Value = 0.4073525654278392
$ws.Range("K4").Value = 0.3819234752393754
$ws.Range("M4").Value = 0.1967081943831701
$ws.Range("N4").Value = 1.191776646103328
$ws.Range("O4").Value = 1.398995391753331
$ws.Range("B5").Value = 0.07219678129175122
$ws.Range("D5").Value = 0.0155677913961938
$ws.Range("E5").Value = 0.1505075918706282
$ws.Range("F5").Value = 0.4332198076258464
$ws.Range("G5").Value = 0.2820811800761334
$ws.Range("H5").Value = 0.4554199336845386
$ws.Range("I5").Value = 0.4088950026180918
$ws.Range("K5").Value = 0.3674088733702092
$ws.Range("M5").Value = 0.1898080765762202
$ws.Range("N5").Value = 1.194729363997673
$ws.Range("O5").Value = 1.398965926131709
$ws.Range("B6").Value = 0.07177052594460065
$ws.Range("D6").Value = 0.01547016302769322
$ws.Range("E6").Value = 0.1501690901071484
$ws.Range("F6").Value = 0.4330791865589276
$ws.Range("G6").Value = 0.2819883664976643
$ws.Range("H6").Value = 0.4555235830923579
$ws.Range("I6").Value = 0.409154180163501
$ws.Range("K6").Value = 0.3649960945885482
$ws.Range("M6").Value = 0.1886628097377425
$ws.Range("N6").Value = 1.195226030645575
$ws.Range("O6").Value = 1.398975248486224
$ws.Range("B7").Value = 0.07472934361749139
$ws.Range("D7").Value = 0.01614711686720938
$ws.Range("E7").Value = 0.1525304870780673
$ws.Range("F7").Value = 0.4340822165965648
$ws.Range("G7").Value = 0.2826532218428497
$ws.Range("H7").Value = 0.4548235011719015
$ws.Range("I7").Value = 0.407373162336421
$ws.Range("K7").Value = 0.3817279037960475
$ws.Range("M7").Value = 0.1966151040231949
$ws.Range("N7").Value = 1.191816040492643
$ws.Range("O7").Value = 1.398994041347393
$ws.Range("B8").Value = 0.08778789520354735
$ws.Range("D8").Value = 0.01911696515944072
$ws.Range("E8").Value = 0.1632448324147546
$ws.Range("F8").Value = 0.4391756377033502
$ws.Range("G8").Value = 0.2860980879263053
$ws.Range("H8").Value = 0.4522129751634267
$ws.Range("I8").Value = 0.3999582786725622
$ws.Range("K8").Value = 0.4551719422474321
$ws.Range("M8").Value = 0.2317778244700079
$ws.Range("N8").Value = 1.177693959836386
$ws.Range("O8").Value = 1.401143941285028
$ws.Range("B9").Value = 0.1133893493917384
$ws.Range("D9").Value = 0.02487629220360077
$ws.Range("E9").Value = 0.1853270435986403
$ws.Range("F9").Value = 0.4515630959664492
$ws.Range("G9").Value = 0.2946963278397305
$ws.Range("H9").Value = 0.448812592736104
$ws.Range("I9").Value = 0.3870146165419186
$ws.Range("K9").Value = 0.5977443645022618
$ws.Range("M9").Value = 0.3009848423491874
$ws.Range("N9").Value = 1.153322261599008
$ws.Range("O9").Value = 1.412798404550898
$ws.Range("B10").Value = 0.1321853186727964
$ws.Range("D10").Value = 0.02906702768230218
$ws.Range("E10").Value = 0.2022199582146271
$ws.Range("F10").Value = 0.4621300795045258
$ws.Range("G10").Value = 0.3021442773333689
$ws.Range("H10").Value = 0.4473609852521179
$ws.Range("I10").Value = 0.3784778456581881
$ws.Range("K10").Value = 0.7015751557678414
$ws.Range("M10").Value = 0.3519902705051976
$ws.Range("N10").Value = 1.137436831292383
$ws.Range("O10").Value = 1.425909114327851
$ws.Range("B11").Value = 0.1407312580395654
$ws.Range("D11").Value = 0.03096442121005794
$ws.Range("E11").Value = 0.2100537993874454
$ws.Range("F11").Value = 0.4672563817015458
$ws.Range("G11").Value = 0.3057797333196532
$ws.Range("H11").Value = 0.4469276629105536
$ws.Range("I11").Value = 0.3748053794811526
$ws.Range("K11").Value = 0.7486049814072828
$ws.Range("M11").Value = 0.3752298055986429
$ws.Range("N11").Value = 1.130647941292324
$ws.Range("O11").Value = 1.432864212931946
$ws.Range("B12").Value = 0.143966542377683
$ws.Range("D12").Value = 0.03168159093229406
$ws.Range("E12").Value = 0.2130419658939644
$ws.Range("F12").Value = 0.469243523433903
$ws.Range("G12").Value = 0.307192059980224
$ws.Range("H12").Value = 0.446796195768286
$ws.Range("I12").Value = 0.3734450468079693
$ws.Range("K12").Value = 0.766383958916947
$ws.Range("M12").Value = 0.3840353063022093
$ws.Range("N12").Value = 1.128140018696861
$ws.Range("O12").Value = 1.435640607334591
$ws.Range("B13").Value = 0.1432698088317039
$ws.Range("D13").Value = 0.03152719546302052
$ws.Range("E13").Value = 0.2123974440531669
$ws.Range("F13").Value = 0.4688135144594554
$ws.Range("G13").Value = 0.3068863028226758
$ws.Range("H13").Value = 0.4468230590315869
$ws.Range("I13").Value = 0.3737366686894443
$ws.Range("K13").Value = 0.7625562962830656
$ws.Range("M13").Value = 0.3821386550455372
$ws.Range("N13").Value = 1.12867734835023
$ws.Range("O13").Value = 1.435036314399468
$ws.Range("B14").Value = 0.1409974457809113
$ws.Range("D14").Value = 0.03102345011329533
$ws.Range("E14").Value = 0.2102992023245491
$ws.Range("F14").Value = 0.4674189446812349
$ws.Range("G14").Value = 0.3058952111238398
$ws.Range("H14").Value = 0.446916193353772
$ws.Range("I14").Value = 0.3746928557329721
$ws.Range("K14").Value = 0.7500682792692999
$ws.Range("M14").Value = 0.3759541360232816
$ws.Range("N14").Value = 1.130440352863715
$ws.Range("O14").Value = 1.433089768816473
$ws.Range("B15").Value = 0.1396054364501396
$ws.Range("D15").Value = 0.03071471691801975
$ws.Range("E15").Value = 0.2090167958023699
$ws.Range("F15").Value = 0.4665707109460158
$ws.Range("G15").Value = 0.3052927852746592
$ws.Range("H15").Value = 0.4469774885699849
$ws.Range("I15").Value = 0.3752825009116192
$ws.Range("K15").Value = 0.7424150472880058
$ws.Range("M15").Value = 0.3721666124270939
$ws.Range("N15").Value = 1.13152843251212
$ws.Range("O15").Value = 1.431916035002672
$ws.Range("B16").Value = 0.1316267116487353
$ws.Range("D16").Value = 0.02894284427193128
$ws.Range("E16").Value = 0.2017110179478721
$ws.Range("F16").Value = 0.4618014896158442
$ws.Range("G16").Value = 0.3019116767186887
$ws.Range("H16").Value = 0.4473938701876392
$ws.Range("I16").Value = 0.3787220962274205
$ws.Range("K16").Value = 0.6984974726831865
$ws.Range("M16").Value = 0.350472245198354
$ws.Range("N16").Value = 1.137889302061126
$ws.Range("O16").Value = 1.425474535943948
$ws.Range("B17").Value = 0.1267307149044541
$ws.Range("D17").Value = 0.02785352685857845
$ws.Range("E17").Value = 0.197267515702805
$ws.Range("F17").Value = 0.4589575172266223
$ws.Range("G17").Value = 0.2999008920608475
$ws.Range("H17").Value = 0.4477074361989821
$ws.Range("I17").Value = 0.3808862179996746
$ws.Range("K17").Value = 0.6715027047825117
$ws.Range("M17").Value = 0.3371728260619804
$ws.Range("N17").Value = 1.141903517822698
$ws.Range("O17").Value = 1.421776803166296
$ws.Range("B18").Value = 0.1239142632876877
$ws.Range("D18").Value = 0.02722613615227942
$ws.Range("E18").Value = 0.1947257587413489
$ws.Range("F18").Value = 0.4573517961086893
$ws.Range("G18").Value = 0.2987676179962335
$ws.Range("H18").Value = 0.4479091607381918
$ws.Range("I18").Value = 0.3821508186282312
$ws.Range("K18").Value = 0.655956947948539
$ws.Range("M18").Value = 0.3295268330020136
$ws.Range("N18").Value = 1.144253566831047
$ws.Range("O18").Value = 1.419743231178245
$ws.Range("B19").Value = 0.1229605990930054
$ws.Range("D19").Value = 0.02701356863880022
$ws.Range("E19").Value = 0.1938675664082226
$ws.Range("F19").Value = 0.4568132886538052
$ws.Range("G19").Value = 0.2983879057923957
$ws.Range("H19").Value = 0.4479811321562721
$ws.Range("I19").Value = 0.3825824003144866
$ws.Range("K19").Value = 0.6506901739721798
$ws.Range("M19").Value = 0.3269386324722348
$ws.Range("N19").Value = 1.145056326264751
$ws.Range("O19").Value = 1.419070712610619
$ws.Range("B20").Value = 0.1272519456364023
$ws.Range("D20").Value = 0.02796957420338231
$ws.Range("E20").Value = 0.1977390807313313
$ws.Range("F20").Value = 0.4592571521710269
$ws.Range("G20").Value = 0.3001125338751081
$ws.Range("H20").Value = 0.4476718450458748
$ws.Range("I20").Value = 0.3806537885087131
$ws.Range("K20").Value = 0.6743783256013103
$ws.Range("M20").Value = 0.3385882124817954
$ws.Range("N20").Value = 1.141471936094604
$ws.Range("O20").Value = 1.422160779352311
$ws.Range("B21").Value = 0.1416649196067965
$ws.Range("D21").Value = 0.03117144879421119
$ws.Range("E21").Value = 0.2109149172379716
$ws.Range("F21").Value = 0.4678273172906486
$ws.Range("G21").Value = 0.3061853503083967
$ws.Range("H21").Value = 0.4468879523584803
$ws.Range("I21").Value = 0.3744111765893914
$ws.Range("K21").Value = 0.7537371397550885
$ws.Range("M21").Value = 0.3777705391537793
$ws.Range("N21").Value = 1.129920809272249
$ws.Range("O21").Value = 1.433657643795783
$ws.Range("B22").Value = 0.1510794063331389
$ws.Range("D22").Value = 0.0332562767566742
$ws.Range("E22").Value = 0.2196524414860832
$ws.Range("F22").Value = 0.4736960868225637
$ws.Range("G22").Value = 0.3103621723188184
$ws.Range("H22").Value = 0.446565772258495
$ws.Range("I22").Value = 0.3705081824350493
$ws.Range("K22").Value = 0.8054262575557232
$ws.Range("M22").Value = 0.4034086701891795
$ws.Range("N22").Value = 1.122737998733818
$ws.Range("O22").Value = 1.44200304733863
$ws.Range("B23").Value = 0.14605527449352
$ws.Range("D23").Value = 0.0321442897671318
$ws.Range("E23").Value = 0.2149774296065701
$ws.Range("F23").Value = 0.4705393213782685
$ws.Range("G23").Value = 0.3081138727396109
$ws.Range("H23").Value = 0.4467203356091716
$ws.Range("I23").Value = 0.3725750910773256
$ws.Range("K23").Value = 0.7778552650817687
$ws.Range("M23").Value = 0.3897223905768499
$ws.Range("N23").Value = 1.126538071119796
$ws.Range("O23").Value = 1.43747281397134
$ws.Range("B24").Value = 0.1270163024308744
$ws.Range("D24").Value = 0.02791711270801045
$ws.Range("E24").Value = 0.1975258461056555
$ws.Range("F24").Value = 0.4591215958943664
$ws.Range("G24").Value = 0.300016779744368
$ws.Range("H24").Value = 0.4476878689848576
$ws.Range("I24").Value = 0.3807588062821914
$ws.Range("K24").Value = 0.6730783386719281
$ws.Range("M24").Value = 0.3379483162127599
$ws.Range("N24").Value = 1.141666922669721
$ws.Range("O24").Value = 1.421986896208637
$ws.Range("B25").Value = 0.1064651446248348
$ws.Range("D25").Value = 0.02332528459988481
$ws.Range("E25").Value = 0.1792368676620555
$ws.Range("F25").Value = 0.4479547608345626
$ws.Range("G25").Value = 0.2921722200521799
$ws.Range("H25").Value = 0.4495485963743846
$ws.Range("I25").Value = 0.3903453207122913
$ws.Range("K25").Value = 0.5593334631096241
$ws.Range("M25").Value = 0.2822349235628607
$ws.Range("N25").Value = 1.159560322235194
$ws.Range("O25").Value = 1.408847821499421
